# Applies the commit: "delete unneeded cells, add simple version for current usage"
# - Inserts a new row 5 containing a new compound (arbr144 / BrC1=CC=CC=C1) marked
#   "untested", with no feature values (D:M left blank).
# - All existing data rows 5-15 shift down to rows 6-16 (content unchanged).
# - Updates the active selection to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 5; this shifts rows 5:15 down to 6:16
# and extends the sheet dimension automatically.
$ws.Rows("5:5").Insert()

# Populate the new row 5 with the new compound entry (only id/smiles/type set,
# feature columns D:M intentionally left empty).
$ws.Range("A5").Value = "arbr144"
$ws.Range("B5").Value = "BrC1=CC=CC=C1"
$ws.Range("C5").Value = "untested"

# Update the sheet's current selection.
$ws.Range("E11").Select()
